$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 517.6224248938721
$ws.Range("D2").Value = 118.5645605878984
$ws.Range("F2").Value = 442
$ws.Range("H2").Value = 552

$ws.Range("C3").Value = 37.34575783424336
$ws.Range("D3").Value = 6.482859032778124
$ws.Range("F3").Value = 32.67
$ws.Range("G3").Value = 37.55
$ws.Range("H3").Value = 41.73

$ws.Range("C4").Value = 2.064765594994014
$ws.Range("D4").Value = 2.53455505977305
$ws.Range("F4").Value = 0.68
$ws.Range("G4").Value = 1.33
$ws.Range("H4").Value = 2.54

$ws.Range("C5").Value = 322.6969035089136
$ws.Range("D5").Value = 8.494172126422965
$ws.Range("F5").Value = 317.8
$ws.Range("G5").Value = 323.21
$ws.Range("H5").Value = 328.4

$ws.Range("C6").Value = 23.74074421510393
$ws.Range("D6").Value = 3.711613999081814
$ws.Range("F6").Value = 21.05
$ws.Range("G6").Value = 23.36
$ws.Range("H6").Value = 26.31

$ws.Range("C7").Value = -75.11641989524432
$ws.Range("D7").Value = 22.05405367114712
$ws.Range("F7").Value = -91
$ws.Range("G7").Value = -71

$ws.Range("C8").Value = 7.9713349218551
$ws.Range("D8").Value = 6.526893324470904
$ws.Range("F8").Value = 8
$ws.Range("H8").Value = 11.2

$ws.Range("C9").Value = 9.112710158057105
$ws.Range("D9").Value = 1.601501744710999

$ws.Range("C10").Value = 867.8228457262176
$ws.Range("D10").Value = 0.4610641877304697

$ws.Range("C11").Value = 0.4718851948034393
$ws.Range("D11").Value = 0.5333773730009826

$ws.Range("C12").Value = 22.74755277862356
$ws.Range("D12").Value = 12.2962608807795

$ws.Range("C13").Value = 0.6716223689530927
$ws.Range("D13").Value = 0.7499756835717255

$ws.Range("C14").Value = 1.831174106852896
$ws.Range("D14").Value = 1.669175369014386

$ws.Range("C15").Value = 92.3764198952438
$ws.Range("D15").Value = 22.05405367114713
$ws.Range("G15").Value = 88.25999999999999
$ws.Range("H15").Value = 108.26

$ws.Range("C16").Value = -84.54313662807775
$ws.Range("D16").Value = 19.91707498201135
$ws.Range("F16").Value = -100.4668316388797
$ws.Range("G16").Value = -82.21081852649533
$ws.Range("H16").Value = -68.5175485570292

$ws.Range("C17").Value = -76.57180170622264
$ws.Range("D17").Value = 24.27898143112502
$ws.Range("F17").Value = -91.2778545523916
$ws.Range("G17").Value = -70.79009749652566
$ws.Range("H17").Value = -57.75746206410165
